# Bugs.xlsx update:
#  - 4 previously-unresolved bugs (D19 "in progress", D20 "iSSUE",
#    D26 "in progress") now have a real resolution date (2019-01-31), and
#    the last row (row 35/36, bug "AD") also gets that resolution date
#    in D35/D36.
#  - Once those are the only cells referencing the shared strings
#    "in progress" / "iSSUE", those strings fall out of the shared-string
#    pool automatically, which re-numbers every other shared-string-backed
#    cell in the sheet (purely mechanical - no extra edits needed for that).
#  - The view was scrolled so row 10 is back at the top of the window.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$resolvedDate = Get-Date -Year 2019 -Month 1 -Day 31 -Hour 0 -Minute 0 -Second 0

# D29 already carries the "short date" cell style used throughout the
# Resolved column - copy its formatting (number format) onto the cells we
# are turning into real dates before writing the date value into them, so
# they pick up the existing date style instead of Excel minting a new one.
$dateStyleSource = $ws.Range("D29")

function Set-ResolvedDate([string]$addr) {
    $dateStyleSource.Copy()
    $ws.Range($addr).PasteSpecial(-4122)   # xlPasteFormats
    $ws.Range($addr).Value = $resolvedDate
}

Set-ResolvedDate "D19"
Set-ResolvedDate "D20"
Set-ResolvedDate "D26"
Set-ResolvedDate "D35"
Set-ResolvedDate "D36"

# Scroll the sheet view back up so row 10 is the top visible row (was 16).
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
